# Updated symbol list on Thu Dec 29 22:37:27 UTC 2022 with GitHub Actions
# Applies the Price (column D) and Volume(1h) (column E) updates described
# by the commit diff to the crypto listing worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (Price) updates -------------------------------------------
# These cells hold text-formatted numbers, so a leading apostrophe is used
# to force Excel to store them as text (preserving formatting such as
# trailing zeros) instead of converting them to numeric values.
$ws.Cells.Item(2, 4).Value  = "'245.23"
$ws.Cells.Item(4, 4).Value  = "'5.280"
$ws.Cells.Item(6, 4).Value  = "'6.459"
$ws.Cells.Item(7, 4).Value  = "'3.128"
$ws.Cells.Item(9, 4).Value  = "'0.8504"
$ws.Cells.Item(10, 4).Value = "'0.1355"
$ws.Cells.Item(11, 4).Value = "'0.06931"
$ws.Cells.Item(12, 4).Value = "'0.03129"
$ws.Cells.Item(13, 4).Value = "'0.02933"
$ws.Cells.Item(15, 4).Value = "'3.744"
$ws.Cells.Item(17, 4).Value = "'0.04683"
$ws.Cells.Item(18, 4).Value = "'0.0005968"
$ws.Cells.Item(19, 4).Value = "'0.006201"
$ws.Cells.Item(21, 4).Value = "'0.004615"
$ws.Cells.Item(22, 4).Value = "'0.00006900"
$ws.Cells.Item(23, 4).Value = "'3.507"
$ws.Cells.Item(25, 4).Value = "'0.3193"
$ws.Cells.Item(26, 4).Value = "'0.1319"
$ws.Cells.Item(28, 4).Value = "'0.0002331"
$ws.Cells.Item(40, 4).Value = "'0.03633"
$ws.Cells.Item(41, 4).Value = "'0.006230"
$ws.Cells.Item(42, 4).Value = "'0.1054"
$ws.Cells.Item(43, 4).Value = "'0.002770"
$ws.Cells.Item(44, 4).Value = "'0.008416"
$ws.Cells.Item(45, 4).Value = "'0.00005250"
$ws.Cells.Item(47, 4).Value = "'0.3699"
$ws.Cells.Item(48, 4).Value = "'0.002357"

# --- Column E (Volume(1h)) updates ---------------------------------------
$ws.Cells.Item(18, 5).Value = "17OneONE"
$ws.Cells.Item(41, 5).Value = "40KickTokenKICK"
$ws.Cells.Item(44, 5).Value = "43LocalTradersLCTBestin24h"
$ws.Cells.Item(47, 5).Value = "46CoinbaseStockTokenCOINWorstin24h"
